$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.519.91"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.916.57"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4795"
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2889"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06729"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.27"
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.19"
$ws.Range("E11").Value = "  +4.60%  "
$ws.Range("D12").Value = "1.906.78"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07574"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.252"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6682"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "301.59"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "30.505.84"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007576"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.557"
$ws.Range("E21").Value = "  +5.88%  "
$ws.Range("D22").Value = "2.161.99"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.483"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.49"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.115"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.172"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.024"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04996"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7368"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02052"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9989"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.723"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.678"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.78"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.024"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4453"
$ws.Range("E42").Value = "  +4.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8659"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.912"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.44"
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.79"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.288"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.259"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1233"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2531"
$ws.Range("E51").Value = "  +2.05%  "
